$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Range("D5:D13").ClearContents()
$ws.Range("D5:D13").Interior.Color = 15983311
